$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new rows (34-36) matching existing pattern rows
$rows = @(
    @(10005, 110033),
    @(10005, 110034),
    @(10005, 110035)
)

$r = 34
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $r++
}

# Update selection/view state
$ws.Range("A37:A1048576").EntireRow.Select()
